$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 / 14 swap: Coin name (B) and Link (C) ---
$ws.Range("B13").Value = "Avalanche"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"

# --- Price (D) updates (force text format so values like "513.25" or
#     "0.340" are not reinterpreted/rounded as numbers) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.540.19"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.332.66"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.25"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.25"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.101"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.63"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.745.79"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.514.04"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.332.69"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.47"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "325.55"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.84"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.73"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.95"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.68"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.39"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.94"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.888"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "154.85"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.58"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.53"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0928"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.560"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.20"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0215"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.17"

# --- Volume(1h) (E) updates ---
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +11.70%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +4.67%  "
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("E38").Value = "  +12.62%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  +5.42%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  +1.79%  "
